# Update the cantonal tax rate ("SteuerfussKanton", column I) for canton SO
# rows 86-105 on the "Staatssteuer" sheet from 100 to 104.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Staatssteuer")
$ws.Activate()

for ($r = 86; $r -le 105; $r++) {
    $ws.Cells.Item($r, 9).Value = 104
}

# Reposition the view / selection to match where the edit was made.
$ws.Application.Goto($ws.Range("A84"), $true)
$ws.Range("I105").Select()
